$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.9864338597395088, 1.245588010707412, 0.40536836597547543, 3.053132205246711)"
$ws.Range("C2").Value = "JSU(-1.4100491586375785, 1.1222653704126113, 2.2836597953585382, 4.226987444184161)"
$ws.Range("D2").Value = "JSU(-0.8550901958369866, 1.0256934614109272, 0.6185473699817055, 2.5104263931726782)"
$ws.Range("E2").Value = "NIG(1.5674637610467492, 1.219358000981356, 3.3502065664215674, 6.015672316119565)"
